# Update Simulation_Manifest rows 138-169 (s137..s168) per the 2021-02-13 commit.
# For each of these rows, column B (#file_name), D (#feedback_1_x), E (#feedback_1_y)
# and H (#feedback_1_theta) get reshuffled/new values; columns A, C, F, G, I are unchanged.
# D/E/H are forced to text ("@") before assignment so the numeric-looking strings are
# stored as text (matching the source workbook's inlineStr/text convention) rather than
# being auto-coerced into numbers by Excel.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{Row=138; B="s137_e101_24019-32_1_5.jpeg"; D="739"; E="1443"; H="71"}
    @{Row=139; B="s138_e103_24019-32_1_7.jpeg"; D="909"; E="666"; H="143"}
    @{Row=140; B="s139_e98_24019-32_1_2.jpeg"; D="1497"; E="1365"; H="34"}
    @{Row=141; B="s140_e113_24019-32_3_1.jpeg"; D="542"; E="652"; H="37"}
    @{Row=142; B="s141_e97_24019-32_1_1.jpeg"; D="1406"; E="568"; H="44"}
    @{Row=143; B="s142_e100_24019-32_1_4.jpeg"; D="2082"; E="1902"; H="122"}
    @{Row=144; B="s143_e108_24019-32_2_4.jpeg"; D="238"; E="1517"; H="136"}
    @{Row=145; B="s144_e109_24019-32_2_3.jpeg"; D="1815"; E="1228"; H="60"}
    @{Row=146; B="s145_e119_24019-32_3_7.jpeg"; D="2087"; E="358"; H="146"}
    @{Row=147; B="s146_e118_24019-32_3_6.jpeg"; D="2493"; E="1233"; H="7"}
    @{Row=148; B="s147_e99_24019-32_1_3.jpeg"; D="2061"; E="249"; H="75"}
    @{Row=149; B="s148_e112_24019-32_2_0.jpeg"; D="1075"; E="679"; H="24"}
    @{Row=150; B="s149_e104_24019-32_1_8.jpeg"; D="379"; E="93"; H="73"}
    @{Row=151; B="s150_e122_24019-32_4_6.jpeg"; D="2416"; E="1811"; H="144"}
    @{Row=152; B="s151_e105_24019-32_2_7.jpeg"; D="1928"; E="217"; H="111"}
    @{Row=153; B="s152_e106_24019-32_2_6.jpeg"; D="785"; E="1165"; H="48"}
    @{Row=154; B="s153_e140_24019-32_2_0.jpeg"; D="221"; E="910"; H="108"}
    @{Row=155; B="s154_e149_24019-32_4_7.jpeg"; D="394"; E="1266"; H="51"}
    @{Row=156; B="s155_e147_24019-32_3_7.jpeg"; D="659"; E="539"; H="163"}
    @{Row=157; B="s156_e151_24019-32_4_5.jpeg"; D="525"; E="945"; H="29"}
    @{Row=158; B="s157_e150_24019-32_4_6.jpeg"; D="297"; E="1371"; H="160"}
    @{Row=159; B="s158_e133_24019-32_2_7.jpeg"; D="829"; E="1217"; H="152"}
    @{Row=160; B="s159_e132_24019-32_1_8.jpeg"; D="1584"; E="714"; H="24"}
    @{Row=161; B="s160_e126_24019-32_1_2.jpeg"; D="1506"; E="278"; H="128"}
    @{Row=162; B="s161_e129_24019-32_1_5.jpeg"; D="1361"; E="1457"; H="15"}
    @{Row=163; B="s162_e127_24019-32_1_3.jpeg"; D="1997"; E="1210"; H="111"}
    @{Row=164; B="s163_e136_24019-32_2_4.jpeg"; D="1757"; E="821"; H="180"}
    @{Row=165; B="s164_e128_24019-32_1_4.jpeg"; D="1613"; E="730"; H="112"}
    @{Row=166; B="s165_e125_24019-32_1_1.jpeg"; D="1227"; E="1474"; H="132"}
    @{Row=167; B="s166_e142_24019-32_3_2.jpeg"; D="564"; E="667"; H="165"}
    @{Row=168; B="s167_e144_24019-32_3_4.jpeg"; D="1093"; E="350"; H="42"}
    @{Row=169; B="s168_e143_24019-32_3_3.jpeg"; D="1655"; E="251"; H="103"}
)

foreach ($u in $updates) {
    $r = $u.Row
    $ws.Range("B$r").Value = $u.B
    $ws.Range("D$r").NumberFormat = "@"
    $ws.Range("D$r").Value = $u.D
    $ws.Range("E$r").NumberFormat = "@"
    $ws.Range("E$r").Value = $u.E
    $ws.Range("H$r").NumberFormat = "@"
    $ws.Range("H$r").Value = $u.H
}
